$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44495
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8500
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 340

# Row 3
$ws.Range("D3").Value = 44165
$ws.Range("J3").Value = 38
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = 8263
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 331

# Row 4
$ws.Range("D4").Value = 44489

# Row 5
$ws.Range("D5").Value = 44167
$ws.Range("J5").Value = 60
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("P5").Value = 340

# Row 6
$ws.Range("D6").Value = 44488
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 8000
$ws.Range("M6").Value = 8500
$ws.Range("P6").Value = 340

# Row 7
$ws.Range("D7").Value = 44161
$ws.Range("J7").Value = 53
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6764
$ws.Range("O7").Value = "Región de O'Higgins"
$ws.Range("P7").Value = 271

# Row 8
$ws.Range("D8").Value = 44162
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7562
$ws.Range("O8").Value = "Región de O'Higgins"
$ws.Range("P8").Value = 302

# Row 9
$ws.Range("D9").Value = 44159
$ws.Range("J9").Value = 42
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6738
$ws.Range("P9").Value = 270

# Row 10
$ws.Range("D10").Value = 44160
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 6500
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 6688
$ws.Range("O10").Value = "Región de O'Higgins"
$ws.Range("P10").Value = 268

# Row 11
$ws.Range("D11").Value = 44466
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11500
$ws.Range("O11").Value = "Región de O'Higgins"
$ws.Range("P11").Value = 460

# Row 12
$ws.Range("D12").Value = 44476
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 7500
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7750
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 310

# Row 13
$ws.Range("D13").Value = 44482
$ws.Range("J13").Value = 120

# Row 14
$ws.Range("D14").Value = 44166
$ws.Range("J14").Value = 56
$ws.Range("K14").Value = 7500
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7804
$ws.Range("O14").Value = "Región de O'Higgins"
$ws.Range("P14").Value = 312

# Row 15
$ws.Range("D15").Value = 44487
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 320

# Row 16
$ws.Range("D16").Value = 44487
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 9000
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 360

# Row 17
$ws.Range("D17").Value = 44473
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 9500
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 9750
$ws.Range("P17").Value = 390

# Row 18
$ws.Range("D18").Value = 44484
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 8500
$ws.Range("M18").Value = 8750
$ws.Range("P18").Value = 350

# Row 19
$ws.Range("D19").Value = 44491
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8500
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 340
